# Updates cryptos price/volume cells (columns D and E) to match
# the refreshed market data, keeping every value as literal text
# (values like '1.00' or '66.508.58' must not become numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''66.508.58'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -1.10%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.448.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -0.63%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  -0.03%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''580.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -2.12%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''175.34'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -1.55%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  -0.02%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  +1.80%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''3.446.76'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -0.74%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = '''  -2.52%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''6.83'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -3.30%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '''  -2.77%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''4.042.72'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -0.74%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''30.83'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -3.43%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''0.131'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -3.40%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''66.540.44'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -1.14%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = '''  -2.79%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''3.447.17'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -0.55%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = '''  -3.72%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = '''  -3.17%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''376.17'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -2.94%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''7.68'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -2.01%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  +0.13%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = '''  -0.17%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  -2.83%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  -1.37%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''0.0000117'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -2.91%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  -4.51%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = '''  -1.89%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +0.12%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''5.83'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -5.23%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''23.85'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +1.67%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''1.98'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -3.51%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  -5.49%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -0.06%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  -4.56%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  -4.85%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''159.59'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -2.56%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''0.877'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +0.86%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''27.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +3.48%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''1.78'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -4.94%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  -3.58%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''6.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -5.37%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''4.45'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -3.37%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''2.693.96'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -4.60%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.0693'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -3.69%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''25.25'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -4.58%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''40.20'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -3.46%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  -1.20%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''320.43'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -4.48%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '''  -3.20%  '
$ws.Range("E51").Style = "Normal"
